$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "This research show us that the AI will learn better policy using Human decision.`r"
